# Fix formatting of scraped floating point numbers (and a few names that
# used a comma as a separator) in the "Importe" column and in the
# "Razon social"/"Nombre Fantasia" columns.
#
# The original values used the Argentine/Spanish number format, e.g.
#   "2.080,00"   (period = thousands separator, comma = decimal separator)
# which gets normalised to the "plain" floating point text form:
#   "2080.00"
#
# A handful of provider-name cells also used a comma as a separator between
# two co-contractors, e.g.
#   "FERNANDEZ MARIO H, GALLICET OSCAR M"
# which gets normalised to use a period instead:
#   "FERNANDEZ MARIO H. GALLICET OSCAR M"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -eq $null) { continue }
        if (-not ($val -is [string])) { continue }
        if ($val.IndexOf(",") -lt 0) { continue }

        $isNumericLooking = $val -match '^[0-9\.]*,[0-9]+$'

        if ($isNumericLooking) {
            # Numeric-looking text: "1.234,56" -> "1234.56"
            $newVal = $val.Replace(".", "").Replace(",", ".")
        } else {
            # Free-text (provider names): replace separator comma with a period
            $newVal = $val.Replace(",", ".")
        }

        if ($newVal -eq $val) { continue }

        if ($isNumericLooking) {
            # Prefix with an apostrophe so the numeric-looking text is kept
            # as a literal string instead of being parsed into a Number;
            # restore the original (default) style afterwards so no visible
            # formatting changes stick.
            $origStyle = $cell.Style
            $cell.Value2 = "'" + $newVal
            $cell.Style = $origStyle
        } else {
            $cell.Value2 = $newVal
        }
    }
}
